$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Kitl"
$ws.Cells.Item(2,3).Value = "Kit"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 111.3936496666667
$ws.Cells.Item(2,8).Value = 334.180949
$ws.Cells.Item(2,9).Value = 0.728503147516836
$ws.Cells.Item(2,10).Value = 0.728503147516836
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 10.20338166666667
$ws.Cells.Item(2,14).Value = 30.610145
$ws.Cells.Item(2,15).Value = 0.9872556892395703
$ws.Cells.Item(2,16).Value = 0.9872556892395703
$ws.Cells.Item(2,17).Value = 1136.591922791956
$ws.Cells.Item(2,18).Value = 10229.32730512761
$ws.Cells.Item(2,19).Value = 0.7192188770149303
$ws.Cells.Item(2,20).Value = 0.7192188770149303

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Kitl"
$ws.Cells.Item(3,3).Value = "Kit"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 111.3936496666667
$ws.Cells.Item(3,8).Value = 334.180949
$ws.Cells.Item(3,9).Value = 0.728503147516836
$ws.Cells.Item(3,10).Value = 0.728503147516836
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.1317136666666667
$ws.Cells.Item(3,14).Value = 0.395141
$ws.Cells.Item(3,15).Value = 0.01274431076042969
$ws.Cells.Item(3,16).Value = 0.01274431076042969
$ws.Cells.Item(3,17).Value = 14.67206604097878
$ws.Cells.Item(3,18).Value = 132.048594368809
$ws.Cells.Item(3,19).Value = 0.00928427050190571
$ws.Cells.Item(3,20).Value = 0.00928427050190571

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Kitl"
$ws.Cells.Item(4,3).Value = "Kit"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 17.89779133333333
$ws.Cells.Item(4,8).Value = 53.693374
$ws.Cells.Item(4,9).Value = 0.1170497363085729
$ws.Cells.Item(4,10).Value = 0.1170497363085729
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 10.20338166666667
$ws.Cells.Item(4,14).Value = 30.610145
$ws.Cells.Item(4,15).Value = 0.9872556892395703
$ws.Cells.Item(4,16).Value = 0.9872556892395703
$ws.Cells.Item(4,17).Value = 182.6179959643589
$ws.Cells.Item(4,18).Value = 1643.56196367923
$ws.Cells.Item(4,19).Value = 0.1155580180946301
$ws.Cells.Item(4,20).Value = 0.1155580180946301

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Kitl"
$ws.Cells.Item(5,3).Value = "Kit"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 17.89779133333333
$ws.Cells.Item(5,8).Value = 53.693374
$ws.Cells.Item(5,9).Value = 0.1170497363085729
$ws.Cells.Item(5,10).Value = 0.1170497363085729
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.1317136666666667
$ws.Cells.Item(5,14).Value = 0.395141
$ws.Cells.Item(5,15).Value = 0.01274431076042969
$ws.Cells.Item(5,16).Value = 0.01274431076042969
$ws.Cells.Item(5,17).Value = 2.357383721748223
$ws.Cells.Item(5,18).Value = 21.216453495734
$ws.Cells.Item(5,19).Value = 0.001491718213942803
$ws.Cells.Item(5,20).Value = 0.001491718213942803

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Kitl"
$ws.Cells.Item(6,3).Value = "Kit"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 23.616134
$ws.Cells.Item(6,8).Value = 70.84840200000001
$ws.Cells.Item(6,9).Value = 0.1544471161745911
$ws.Cells.Item(6,10).Value = 0.1544471161745911
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 10.20338166666667
$ws.Cells.Item(6,14).Value = 30.610145
$ws.Cells.Item(6,15).Value = 0.9872556892395703
$ws.Cells.Item(6,16).Value = 0.9872556892395703
$ws.Cells.Item(6,17).Value = 240.9644286931434
$ws.Cells.Item(6,18).Value = 2168.67985823829
$ws.Cells.Item(6,19).Value = 0.1524787941300099
$ws.Cells.Item(6,20).Value = 0.1524787941300099

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Kitl"
$ws.Cells.Item(7,3).Value = "Kit"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 23.616134
$ws.Cells.Item(7,8).Value = 70.84840200000001
$ws.Cells.Item(7,9).Value = 0.1544471161745911
$ws.Cells.Item(7,10).Value = 0.1544471161745911
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.1317136666666667
$ws.Cells.Item(7,14).Value = 0.395141
$ws.Cells.Item(7,15).Value = 0.01274431076042969
$ws.Cells.Item(7,16).Value = 0.01274431076042969
$ws.Cells.Item(7,17).Value = 3.110567601631334
$ws.Cells.Item(7,18).Value = 27.995108414682
$ws.Cells.Item(7,19).Value = 0.001968322044581176
$ws.Cells.Item(7,20).Value = 0.001968322044581176

